# "Added affix types for filters"
#
# Adds a new "affix_type" column (AT) to the Affixes / stat_modifiers sheet:
#   - AT1 gets the new header "affix_type"
#   - AT2:AT59 (every existing data row) gets a default value of 0
#   - the new column is given a width matching its neighbours
#   - the selection is moved onto the freshly added column, like a user
#     would do right after typing the new data in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header -------------------------------------------------------
$ws.Range("AT1").Value = "affix_type"

# --- New data column, defaulting to 0 for every existing row ----------
$ws.Range("AT2:AT59").Value = 0

# --- Column width for the new column (matches the ~17.24 char width) --
$ws.Columns.Item(46).ColumnWidth = 16.28

# --- Keep the view settings that were already on the sheet -----------
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true
$excel.ActiveWindow.DisplayZeros = $true
$excel.ActiveWindow.DisplayFormulas = $false

# --- Move the selection onto the new column, like a user just would --
$ws.Range("AT2:AT59").Select()
